$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 42. This shifts the existing rows 42:64 down to 43:65,
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with a new weekly price record.
$ws.Cells.Item(42, 1).Value = 2
$ws.Cells.Item(42, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 44875
$ws.Cells.Item(42, 5).Value = 4
$ws.Cells.Item(42, 6).Value = 100112032
$ws.Cells.Item(42, 7).Value = "Zapallo italiano"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 500
$ws.Cells.Item(42, 11).Value = 8000
$ws.Cells.Item(42, 12).Value = 9000
$ws.Cells.Item(42, 13).Value = 8500
$ws.Cells.Item(42, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 142
$ws.Cells.Item(42, 17).Value = 60
$ws.Cells.Item(42, 18).Value = "Hortaliza"
